$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.528.19'
$ws.Range('E2').Value = '  +0.79%  '
$ws.Range('D3').Value = '3.755.11'
$ws.Range('E3').Value = '  -0.68%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = "'593.89"
$ws.Range('E5').Value = '  -0.46%  '
$ws.Range('D6').Value = "'166.92"
$ws.Range('E6').Value = '  -1.87%  '
$ws.Range('D7').Value = '3.753.80'
$ws.Range('E7').Value = '  -0.70%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('E9').Value = '  -1.11%  '
$ws.Range('E10').Value = '  -2.80%  '
$ws.Range('E11').Value = '  -1.27%  '
$ws.Range('D12').Value = "'0.449"
$ws.Range('E12').Value = '  -1.17%  '
$ws.Range('E13').Value = '  -7.50%  '
$ws.Range('D14').Value = "'36.05"
$ws.Range('E14').Value = '  -1.77%  '
$ws.Range('D15').Value = '4.384.36'
$ws.Range('E15').Value = '  -0.69%  '
$ws.Range('D16').Value = '3.753.92'
$ws.Range('E16').Value = '  -0.85%  '
$ws.Range('D17').Value = '68.523.91'
$ws.Range('E17').Value = '  +0.92%  '
$ws.Range('D18').Value = "'17.96"
$ws.Range('E18').Value = '  -4.96%  '
$ws.Range('E19').Value = '  +0.66%  '
$ws.Range('E20').Value = '  -3.08%  '
$ws.Range('D21').Value = "'10.74"
$ws.Range('E21').Value = '  +1.17%  '
$ws.Range('D22').Value = "'464.37"
$ws.Range('E22').Value = '  -0.75%  '
$ws.Range('D23').Value = "'0.697"
$ws.Range('E23').Value = '  -3.26%  '
$ws.Range('D24').Value = "'84.10"
$ws.Range('E24').Value = '  +0.45%  '
$ws.Range('E25').Value = '  -3.42%  '
$ws.Range('D26').Value = "'2.17"
$ws.Range('E26').Value = '  -3.36%  '
$ws.Range('D27').Value = "'11.94"
$ws.Range('E27').Value = '  -1.85%  '
$ws.Range('B28').Value = 'Dai'
$ws.Range('C28').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D28').Value = "'1.00"
$ws.Range('E28').Value = '  -0.10%  '
$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D29').Value = "'10.00"
$ws.Range('E29').Value = '  -5.05%  '
$ws.Range('D30').Value = '3.901.33'
$ws.Range('E30').Value = '  -0.72%  '
$ws.Range('E31').Value = '  -5.49%  '
$ws.Range('D32').Value = "'7.31"
$ws.Range('E32').Value = '  -3.90%  '
$ws.Range('D33').Value = "'29.98"
$ws.Range('E33').Value = '  -1.82%  '
$ws.Range('E34').Value = '  -3.32%  '
$ws.Range('E35').Value = '  -0.78%  '
$ws.Range('D36').Value = "'0.997"
$ws.Range('D37').Value = '3.706.95'
$ws.Range('E37').Value = '  -0.82%  '
$ws.Range('E38').Value = '  -4.23%  '
$ws.Range('D39').Value = "'3.39"
$ws.Range('E39').Value = '  -10.21%  '
$ws.Range('E40').Value = '  -0.43%  '
$ws.Range('E41').Value = '  -0.44%  '
$ws.Range('D42').Value = "'5.78"
$ws.Range('E42').Value = '  -1.07%  '
$ws.Range('E44').Value = '  +0.00%  '
$ws.Range('D45').Value = "'44.10"
$ws.Range('E45').Value = '  +9.73%  '
$ws.Range('E46').Value = '  -4.02%  '
$ws.Range('B47').Value = 'OKB'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D47').Value = "'46.85"
$ws.Range('E47').Value = '  +2.65%  '
$ws.Range('B48').Value = 'Stacks'
$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D48').Value = "'1.91"
$ws.Range('E48').Value = '  -2.22%  '
$ws.Range('B49').Value = 'Cosmos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D49').Value = "'8.48"
$ws.Range('E49').Value = '  -2.46%  '
$ws.Range('D50').Value = "'145.41"
$ws.Range('E50').Value = '  +3.07%  '
$ws.Range('D51').Value = "'387.56"
$ws.Range('E51').Value = '  -4.21%  '
